$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data table (rows 2..11), columns A..I
$data = @(
    @("2025-07-28", 2, "BEMOL S/A", "389107", 10525, "BOMBA AUTOMATICA PARA GALAO DE AGUA RECARREGAVEL USB", -113, 1.03, 0.18),
    @("2025-07-29", 2, "BEMOL S/A", "389675", 49177, "KIT SMARTWATCH INOVA PULSEIRA PRETO LISA (FONE+FONTE+CABO)", -69, 1.03, 0.17),
    @("2025-07-29", 2, "BEMOL S/A", "389699", 48696, "FONE DE OUVIDO TIPO-C EJ-105", -76, 1.05, 0.23),
    @("2025-07-30", 2, "BEMOL S/A", "390273", 13958, "CAMPAINHA RESIDENCIAL SEM FIO 32 TOQUES C/PLUG BIVOL", -26, 1.07, 0.27),
    @("2025-07-30", 2, "BEMOL S/A", "390286", 12945, "FONE DE OUVIDO SEM FIO BT BASIKE FON-9856", -72, 1.02, 0.16),
    @("2025-07-31", 7, "BEMOL S/A", "390878", 48587, "SUPORTE PARA NOTEBOOK HMASTON PRETO", -130, 1.27, 1.46),
    @("2025-08-04", 2, "BEMOL S/A", "391921", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1270, 1.05, 0.26),
    @("2025-08-07", 2, "BEMOL S/A", "393760", 14152, "HEADSET GAMER PLAYER PLUS LED 7 CORES 2M DRIVER 50MM PRETO LETRON", -13, 1.08, 0.29),
    @("2025-08-07", 3, "BEMOL S/A", "393791", 13965, "MINI VENTILADOR RECARREGAVEL E PORTATIL, COM PREGADOR MATERIAL ABS E COMPONENTES ELETRONICOS", -22, 1.24, 0.54),
    @("2025-08-08", 2, "BEMOL S/A", "394429", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1270, 1.05, 0.26)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $row++
}
